# This edit reorders the data rows (2-13) of the "Artfynd" sheet: each
# destination row ends up containing the full set of values that used to
# live in a different source row (row 1, the header, is untouched). The
# mapping below was derived by diffing the Id column (A) plus full-row
# content before/after.
#
# destination row -> source row
$rowMap = @{
    2  = 4
    3  = 5
    4  = 7
    5  = 13
    6  = 2
    7  = 3
    8  = 6
    9  = 8
    10 = 9
    11 = 10
    12 = 11
    13 = 12
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$numRows = $used.Rows.Count
$numCols = $used.Columns.Count

# Read the whole used range into a 2D array. COM hands back a 1-based
# array ([1..numRows, 1..numCols]); the array we assign back must match
# the *exact* extent of the target range ([0..numRows-1, 0..numCols-1],
# .NET's native 0-based indexing) or Excel misaligns the write.
$original = $used.Value2

$updated = New-Object 'object[,]' $numRows, $numCols

for ($c = 1; $c -le $numCols; $c++) {
    $updated[0, $c - 1] = $original[1, $c]
}

# Cells holding plain text that happens to look like a date (e.g.
# "2016-07-05") get auto-coerced into a date serial number when pushed
# back through .Value2 as part of a bulk array write. Remember which
# destination cells need to be re-applied as literal text afterwards.
$dateLikePattern = '^\d{4}-\d{2}-\d{2}$'
$textFixups = New-Object System.Collections.Generic.List[object]

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    for ($c = 1; $c -le $numCols; $c++) {
        $val = $original[$srcRow, $c]
        $updated[$destRow - 1, $c - 1] = $val
        if (($val -is [string]) -and ($val -match $dateLikePattern)) {
            [void]$textFixups.Add(@{ Row = $destRow; Col = $c; Value = $val })
        }
    }
}

$used.Value2 = $updated

# Re-apply the date-like strings as literal text so they keep their
# original textual representation instead of becoming date serials.
foreach ($fix in $textFixups) {
    $cell = $ws.Cells.Item($firstRow + $fix.Row - 1, $firstCol + $fix.Col - 1)
    $cell.NumberFormat = "@"
    $cell.Value2 = $fix.Value
    $cell.Style = "Normal"
}
